$d = $word.ActiveDocument

# Update the date paragraph (first paragraph of the document)
$d.Paragraphs.Item(1).Range.Text = "2023-09-10 Sunday"

# Update each answer cell in the single 20x5 table, in row-major order
$t = $d.Tables.Item(1)
$answers = @(
    "33+57=90",
    "80+13=93",
    "53+14=67",
    "67+12=79",
    "84-14=70",
    "75-43=32",
    "31+19=50",
    "68-43=25",
    "97-77=20",
    "79-79=0",
    "36-13=23",
    "9+48=57",
    "47+34=81",
    "31-7=24",
    "20+41=61",
    "61-11=50",
    "60-14=46",
    "54-11=43",
    "74-23=51",
    "73-55=18",
    "58-31=27",
    "68-67=1",
    "53+10=63",
    "44-5=39",
    "64-3=61",
    "66-16=50",
    "11+29=40",
    "6+15=21",
    "55+42=97",
    "27+1=28",
    "67-11=56",
    "23-2=21",
    "94-70=24",
    "70-13=57",
    "11+42=53",
    "8+20=28",
    "14+23=37",
    "21+65=86",
    "72-56=16",
    "77-13=64",
    "25+19=44",
    "86-79=7",
    "52+36=88",
    "15+43=58",
    "38-2=36",
    "21+51=72",
    "23+75=98",
    "76+0=76",
    "54-20=34",
    "73-1=72",
    "33+46=79",
    "15+38=53",
    "79-25=54",
    "40+55=95",
    "6+87=93",
    "18-14=4",
    "75+10=85",
    "27+3=30",
    "35-5=30",
    "70-50=20",
    "99-56=43",
    "70-13=57",
    "24+69=93",
    "92-31=61",
    "47+2=49",
    "72-38=34",
    "91-39=52",
    "7+50=57",
    "24+53=77",
    "82-31=51",
    "40+31=71",
    "44+38=82",
    "79+1=80",
    "13+38=51",
    "39+14=53",
    "53-6=47",
    "66+17=83",
    "33-15=18",
    "74+2=76",
    "5+9=14",
    "34+48=82",
    "86-62=24",
    "85+5=90",
    "71-43=28",
    "71-59=12",
    "79-47=32",
    "61-28=33",
    "40-33=7",
    "37+43=80",
    "39+19=58",
    "29+45=74",
    "17+14=31",
    "53+4=57",
    "3+18=21",
    "0+40=40",
    "77-37=40",
    "32+19=51",
    "29+49=78",
    "95-84=11",
    "24+18=42"
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $answers[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done. Updated" $idx "cells."
